$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text (Binance rates) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldText = $ws1.Range("A1").Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 1.93 = 7079.34 pesos"), "1000 Bs = 1.94 = 7134.13 pesos"
$newText = $newText -replace [regex]::Escape("7079.34 pesos = 1.92 = 965.75 Bs"), "7134.13 pesos = 1.93 = 947.35 Bs"
$ws1.Range("A1").Value2 = $newText

# --- Update "tasas" sheet rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 515.9
$ws2.Range("O10").Value = 3680.5
$ws2.Range("N12").Value = 3690
$ws2.Range("O12").Value = 490
